$wb = $excel.ActiveWorkbook

# DatosCuenta sheet
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokPreProdJuneTres"
$wsCuenta.Range("B2").Value = "SmokeNameJuneTres"
$wsCuenta.Range("C2").Value = 27100123
$wsCuenta.Range("D2").Value = 124

# DatosHogar sheet
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 643

# DatosMotor sheet
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP024"
$wsMotor.Range("B2").Value = "ABC12SSMP024"
$wsMotor.Range("C2").Value = "ZAZ123SSMP024"

# DatosAP sheet
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200125
$wsAP.Select()
$wsAP.Range("H10").Select()
